$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("29-10-2021", "30-10-2021", "31-10-2021", "01-11-2021", "02-11-2021")

$startRow = 303
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $cell = $ws.Cells.Item($row, 1)
    # Force text so day/month-looking strings (e.g. "01-11-2021") are not
    # auto-converted into date serial numbers by Excel's smart input parsing.
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$i]
    # Restore the default (unformatted) style so the cell matches the rest
    # of the column, which carries no explicit style/number format.
    $cell.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = 1694
    $ws.Cells.Item($row, 3).Value = 2114
    $ws.Cells.Item($row, 4).Value = 12227
    $ws.Cells.Item($row, 5).Value = 2134
    $ws.Cells.Item($row, 6).Value = 3696
    $ws.Cells.Item($row, 7).Value = 7100
}
